$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the test-run results for the two test cases (rows 3 & 4) ---
# Row 3 = FCCB_TB_JCX test case seq 1
$msgRow3 = @'
未知异常No frame element found by name or id mainFrame
Build info: version: '3.14.0', revision: 'aacccce0', time: '2018-08-02T20:19:58.91Z'
System info: host: 'DESKTOP-I9S3ABE', ip: '100.100.90.19', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '1.8.0_101'
Driver info: driver.version: unknown
'@
$shotRow3 = @'
C:\Users\liu-b\git\PageObjectModel\screenshots\FCCB_TB_JCX_1_异常截屏_20200701060502602.png
'@

$ws.Range("BY3").Value = "06201890101202000000000171"
$ws.Range("BZ3").Value = $msgRow3
$ws.Range("CA3").Value = $shotRow3

# Row 4 = FCCB_TB_JCX test case seq 2
$msgRow4 = @'
未知弹窗_投保人 李宇轩 的社会统一信用代码错误，请修改！
'@
$shotRow4 = @'
C:\Users\liu-b\git\PageObjectModel\screenshots\FCCB_TB_JCX_2_未知弹窗_20200701062530719.png
'@

$ws.Range("BY4").Value = "06201890101202000000000172"
$ws.Range("BZ4").Value = $msgRow4
$ws.Range("CA4").Value = $shotRow4

# --- Move the selection/view back to L16 (no frozen scroll position) ---
$ws.Range("L16").Select() | Out-Null
